# Update cryptocurrency price/volume data per latest scrape (Thu Jun 15 08:41:37 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.957.06"
$ws.Range("E2").Value = "  -3.76%  "

$ws.Range("D3").Value = "'1.638.03"
$ws.Range("E3").Value = "  -6.11%  "

$ws.Range("D4").Value = "'0.9967"
$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "'236.57"
$ws.Range("E5").Value = "  -4.41%  "

$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").Value = "'0.4725"
$ws.Range("E7").Value = "  -6.20%  "

$ws.Range("D8").Value = "'0.2566"
$ws.Range("E8").Value = "  -6.09%  "

$ws.Range("D9").Value = "'0.06013"
$ws.Range("E9").Value = "  -2.83%  "

$ws.Range("D10").Value = "'0.07046"
$ws.Range("E10").Value = "  -2.83%  "

$ws.Range("D11").Value = "'1.635.64"
$ws.Range("E11").Value = "  -6.16%  "

$ws.Range("E12").Value = "  -2.30%  "

$ws.Range("D13").Value = "'0.6167"
$ws.Range("E13").Value = "  -5.46%  "

$ws.Range("D14").Value = "'4.369"
$ws.Range("E14").Value = "  -5.70%  "

$ws.Range("D15").Value = "'72.76"
$ws.Range("E15").Value = "  -6.16%  "

$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").Value = "'0.9984"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").Value = "'24.952.43"
$ws.Range("E18").Value = "  -3.86%  "

$ws.Range("D19").Value = "'0.000006587"
$ws.Range("E19").Value = "  -3.35%  "

$ws.Range("D20").Value = "'11.16"
$ws.Range("E20").Value = "  -5.81%  "

$ws.Range("D21").Value = "'4.409"
$ws.Range("E21").Value = "  +1.62%  "

$ws.Range("D22").Value = "'1.847.49"
$ws.Range("E22").Value = "  -6.42%  "

$ws.Range("D23").Value = "'8.607"
$ws.Range("E23").Value = "  -0.77%  "

$ws.Range("D24").Value = "'5.276"
$ws.Range("E24").Value = "  -2.36%  "

$ws.Range("D25").Value = "'133.18"
$ws.Range("E25").Value = "  -2.54%  "

$ws.Range("E26").Value = "  -2.68%  "

$ws.Range("D27").Value = "'1.358"
$ws.Range("E27").Value = "  -9.38%  "

$ws.Range("D28").Value = "'102.62"
$ws.Range("E28").Value = "  -3.01%  "

$ws.Range("D29").Value = "'1.659"
$ws.Range("E29").Value = "  -6.47%  "

$ws.Range("D30").Value = "'3.750"
$ws.Range("E30").Value = "  -4.34%  "

$ws.Range("D31").Value = "'0.07724"
$ws.Range("E31").Value = "  -6.26%  "

$ws.Range("D32").Value = "'3.559"
$ws.Range("E32").Value = "  -1.97%  "

$ws.Range("B33").Value = "Frax"
$ws.Range("C33").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D33").Value = "'0.9984"
$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.04315"
$ws.Range("E34").Value = "  -7.67%  "

$ws.Range("D35").Value = "'2.599"

$ws.Range("D36").Value = "'0.9213"
$ws.Range("E36").Value = "  -7.54%  "

$ws.Range("D37").Value = "'0.5821"
$ws.Range("E37").Value = "  -6.04%  "

$ws.Range("E38").Value = "  -6.04%  "

$ws.Range("D39").Value = "'0.01554"
$ws.Range("E39").Value = "  -3.22%  "

$ws.Range("D40").Value = "'0.9980"
$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("D41").Value = "'0.8333"
$ws.Range("E41").Value = "  +10.10%  "

$ws.Range("D42").Value = "'1.800"
$ws.Range("E42").Value = "  -6.13%  "

$ws.Range("D43").Value = "'97.48"
$ws.Range("E43").Value = "  -2.23%  "

$ws.Range("D44").Value = "'0.3717"
$ws.Range("E44").Value = "  -4.38%  "

$ws.Range("D45").Value = "'4.741"
$ws.Range("E45").Value = "  -5.13%  "

$ws.Range("D46").Value = "'0.1104"
$ws.Range("E46").Value = "  -3.57%  "

$ws.Range("D47").Value = "'0.05217"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("D48").Value = "'6.077"
$ws.Range("E48").Value = "  -3.84%  "

$ws.Range("D49").Value = "'29.61"
$ws.Range("E49").Value = "  -3.26%  "

$ws.Range("D50").Value = "'0.9978"
$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("D51").Value = "'0.9975"
$ws.Range("E51").Value = "  -0.49%  "
